$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "314.43")
# are preserved exactly as literal text instead of being parsed as numbers,
# matching the inline string cells in the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.235.59"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "2.349.79"
$ws.Range("E3").Value = "  +6.07%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "314.43"
$ws.Range("E5").Value = "  +6.69%  "
$ws.Range("D6").Value = "109.60"
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +3.90%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("D10").Value = "42.94"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "0.0938"
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("D12").Value = "8.88"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +9.30%  "
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "16.31"
$ws.Range("E15").Value = "  +9.27%  "
$ws.Range("D16").Value = "2.703.45"
$ws.Range("E16").Value = "  +6.09%  "
$ws.Range("D17").Value = "2.344.67"
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "43.211.65"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "75.35"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("E22").Value = "  +14.37%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "254.02"
$ws.Range("E24").Value = "  +11.57%  "
$ws.Range("D25").Value = "9.11"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("D28").Value = "39.28"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").Value = "22.39"
$ws.Range("E30").Value = "  +7.37%  "
$ws.Range("D31").Value = "174.80"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").Value = "0.0928"
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").Value = "6.01"
$ws.Range("E34").Value = "  +9.32%  "
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "4.15"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").Value = "  +3.99%  "
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  +10.80%  "
$ws.Range("D41").Value = "73.03"
$ws.Range("E41").Value = "  +2.71%  "
$ws.Range("D42").Value = "1.48"
$ws.Range("E42").Value = "  +13.62%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "12.85"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("E47").Value = "  +10.56%  "
$ws.Range("D48").Value = "111.15"
$ws.Range("E48").Value = "  +7.73%  "
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("D51").Value = "69.88"
$ws.Range("E51").Value = "  +4.99%  "

# Restore column D formatting to General so no residual text-format style
# is left behind on the cells (keeps styling identical to the original file).
$ws.Range("D2:D51").ClearFormats()
